$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their text formatting so values like "1.00" or "0.999"
# are not silently converted to numbers (losing trailing zeros / dot formatting).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '98.560.83'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.09%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.371.87'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +7.69%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '259.04'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +6.69%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '628.24'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.69%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +22.74%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.392'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.91%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.05%  '

# Row 10
$ws.Range("B10").Value = 'LidoStakedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.377.70'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +7.98%  '

# Row 11
$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.861'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +9.78%  '

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.11%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '98.519.88'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.41%  '

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +6.38%  '

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.49%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.002.69'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +7.69%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.50'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.27%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.366.79'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +7.40%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.57'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.53%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.22'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.73%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '489.33'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.01%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.07'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +6.84%  '

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +9.59%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.40'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.72%  '

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.81%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.99'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.05%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.97'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.22%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.555.25'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +7.62%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.283'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +14.88%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.02%  '

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +8.66%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.136'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +9.04%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.65'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +7.10%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.30%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '28.02'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.22%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.61%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.31'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.41%  '

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.52%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '501.59'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.80%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.46%  '

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.82%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.75'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.97%  '

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.57%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.30'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.90%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.792'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +13.40%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.01%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '159.80'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.14%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.94'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.10%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.846'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +12.60%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.63'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.30%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '45.92'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.81%  '
